$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new report-date headers for columns BD and BE ---
$ws.Range("BD1").Value = 43976
$ws.Range("BE1").Value = 43978

# --- Rows 2-19: new daily figures for columns BD and BE ---
$values = @{
    2  = @(109, 131)
    3  = @(95, 99)
    4  = @(42, 42)
    5  = @(35, 62)
    6  = @(18, 22)
    7  = @(2680, 2982)
    8  = @(23, 27)
    9  = @(835, 882)
    10 = @(0, 0)
    11 = @(13, 13)
    12 = @(1, 1)
    13 = @(9, 15)
    14 = @(13, 14)
    15 = @(3, 3)
    16 = @(13, 15)
    17 = @(41, 48)
    18 = @(45, 52)
    19 = @(214, 232)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Range("BD$row").Value = $pair[0]
    $ws.Range("BE$row").Value = $pair[1]
}

# --- Row 20: totals (sum formulas matching the existing BC20 pattern) ---
$ws.Range("BD20").Formula = "=SUM(BD2:BD19)"
$ws.Range("BE20").Formula = "=SUM(BE2:BE19)"

# --- Update the sheet's view/selection state to match the last user action ---
$excel.ActiveWindow.DisplayGridlines = $true
[void]$ws.Range("BF20").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 51
